$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the section headers: B6 becomes "Failed Banks", B11 becomes "Survived Banks"
$ws.Range("B6").Value = "Failed Banks"
$ws.Range("B11").Value = "Survived Banks"

# Update the active selection to C15
$ws.Range("C15").Select()
